$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.072131
$ws.Range("H2").Value = 18.216393
$ws.Range("I2").Value = 0.003943999267036455
$ws.Range("J2").Value = 0.003943999267036454
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 122.328922
$ws.Range("N2").Value = 366.986766
$ws.Range("O2").Value = 0.9783373008518612
$ws.Range("P2").Value = 0.9783373008518613
$ws.Range("Q2").Value = 742.7972394727819
$ws.Range("R2").Value = 6685.175155255038
$ws.Range("S2").Value = 0.003858561597474164
$ws.Range("T2").Value = 0.003858561597474164

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.072131
$ws.Range("H3").Value = 18.216393
$ws.Range("I3").Value = 0.003943999267036455
$ws.Range("J3").Value = 0.003943999267036454
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.3863573333333334
$ws.Range("N3").Value = 1.159072
$ws.Range("O3").Value = 0.003089929874945324
$ws.Range("P3").Value = 0.003089929874945324
$ws.Range("Q3").Value = 2.346012340810667
$ws.Range("R3").Value = 21.114111067296
$ws.Range("S3").Value = 0.0000121866811619784
$ws.Range("T3").Value = 0.0000121866811619784

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.072131
$ws.Range("H4").Value = 18.216393
$ws.Range("I4").Value = 0.003943999267036455
$ws.Range("J4").Value = 0.003943999267036454
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.322294
$ws.Range("N4").Value = 6.966882000000001
$ws.Range("O4").Value = 0.0185727692731934
$ws.Range("P4").Value = 0.0185727692731934
$ws.Range("Q4").Value = 14.101273388514
$ws.Range("R4").Value = 126.911460496626
$ws.Range("S4").Value = 0.00007325098840031197
$ws.Range("T4").Value = 0.00007325098840031196

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1480.851806666667
$ws.Range("H5").Value = 4442.55542
$ws.Range("I5").Value = 0.9618498744646554
$ws.Range("J5").Value = 0.9618498744646552
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 122.328922
$ws.Range("N5").Value = 366.986766
$ws.Range("O5").Value = 0.9783373008518612
$ws.Range("P5").Value = 0.9783373008518613
$ws.Range("Q5").Value = 181151.0051512857
$ws.Range("R5").Value = 1630359.046361572
$ws.Range("S5").Value = 0.9410136100084525
$ws.Range("T5").Value = 0.9410136100084524

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1480.851806666667
$ws.Range("H6").Value = 4442.55542
$ws.Range("I6").Value = 0.9618498744646554
$ws.Range("J6").Value = 0.9618498744646552
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.3863573333333334
$ws.Range("N6").Value = 1.159072
$ws.Range("O6").Value = 0.003089929874945324
$ws.Range("P6").Value = 0.003089929874945324
$ws.Range("Q6").Value = 572.1379550855822
$ws.Range("R6").Value = 5149.24159577024
$ws.Range("S6").Value = 0.002972048662320748
$ws.Range("T6").Value = 0.002972048662320747

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1480.851806666667
$ws.Range("H7").Value = 4442.55542
$ws.Range("I7").Value = 0.9618498744646554
$ws.Range("J7").Value = 0.9618498744646552
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.322294
$ws.Range("N7").Value = 6.966882000000001
$ws.Range("O7").Value = 0.0185727692731934
$ws.Range("P7").Value = 0.0185727692731934
$ws.Range("Q7").Value = 3438.97326551116
$ws.Range("R7").Value = 30950.75938960044
$ws.Range("S7").Value = 0.01786421579388208
$ws.Range("T7").Value = 0.01786421579388208

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 52.663316
$ws.Range("H8").Value = 157.989948
$ws.Range("I8").Value = 0.03420612626830831
$ws.Range("J8").Value = 0.0342061262683083
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 122.328922
$ws.Range("N8").Value = 366.986766
$ws.Range("O8").Value = 0.9783373008518612
$ws.Range("P8").Value = 0.9783373008518613
$ws.Range("Q8").Value = 6442.246675225352
$ws.Range("R8").Value = 57980.22007702816
$ws.Range("S8").Value = 0.0334651292459347
$ws.Range("T8").Value = 0.0334651292459347

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 52.663316
$ws.Range("H9").Value = 157.989948
$ws.Range("I9").Value = 0.03420612626830831
$ws.Range("J9").Value = 0.0342061262683083
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.3863573333333334
$ws.Range("N9").Value = 1.159072
$ws.Range("O9").Value = 0.003089929874945324
$ws.Range("P9").Value = 0.003089929874945324
$ws.Range("Q9").Value = 20.34685833425067
$ws.Range("R9").Value = 183.121725008256
$ws.Range("S9").Value = 0.0001056945314625978
$ws.Range("T9").Value = 0.0001056945314625978

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 52.663316
$ws.Range("H10").Value = 157.989948
$ws.Range("I10").Value = 0.03420612626830831
$ws.Range("J10").Value = 0.0342061262683083
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.322294
$ws.Range("N10").Value = 6.966882000000001
$ws.Range("O10").Value = 0.0185727692731934
$ws.Range("P10").Value = 0.0185727692731934
$ws.Range("Q10").Value = 122.299702766904
$ws.Range("R10").Value = 1100.697324902136
$ws.Range("S10").Value = 0.0006353024909110103
$ws.Range("T10").Value = 0.0006353024909110102

